$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.702.90'
$ws.Range("E2").Value = '  -5.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.368.71'
$ws.Range("E3").Value = '  -6.58%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.06'
$ws.Range("E5").Value = '  -5.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.48'
$ws.Range("E6").Value = '  -9.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.600'
$ws.Range("E7").Value = '  -4.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.361.33'
$ws.Range("E9").Value = '  -6.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.190'
$ws.Range("E10").Value = '  -11.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.599'
$ws.Range("E11").Value = '  -7.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.09'
$ws.Range("E12").Value = '  -10.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  -10.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.75'
$ws.Range("E14").Value = '  -9.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.899.08'
$ws.Range("E15").Value = '  -6.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '609.76'
$ws.Range("E16").Value = '  -10.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.602.99'
$ws.Range("E17").Value = '  -5.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.369.85'
$ws.Range("E18").Value = '  -6.70%  '
$ws.Range("E19").Value = '  -4.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.77'
$ws.Range("E20").Value = '  -6.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.70'
$ws.Range("E21").Value = '  -8.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.915'
$ws.Range("E22").Value = '  -8.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.03'
$ws.Range("E23").Value = '  -7.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.09'
$ws.Range("E24").Value = '  -3.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.95'
$ws.Range("E25").Value = '  -11.62%  '
$ws.Range("E26").Value = '  -10.10%  '
$ws.Range("E27").Value = '  -8.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.54'
$ws.Range("E28").Value = '  -9.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.87'
$ws.Range("E29").Value = '  -11.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.99'
$ws.Range("E30").Value = '  -9.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.38'
$ws.Range("E31").Value = '  -10.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.86'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.25'
$ws.Range("E33").Value = '  -8.67%  '
$ws.Range("E34").Value = '  -7.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.838.86'
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.93'
$ws.Range("E36").Value = '  -7.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '539.60'
$ws.Range("E37").Value = '  +5.62%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.81'
$ws.Range("E39").Value = '  +39.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  -5.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0729'
$ws.Range("E41").Value = '  -13.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("E42").Value = '  -9.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.354'
$ws.Range("E43").Value = '  -7.67%  '
$ws.Range("E44").Value = '  -6.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.86'
$ws.Range("E45").Value = '  -10.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0420'
$ws.Range("E46").Value = '  -10.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.69'
$ws.Range("E47").Value = '  -12.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.14'
$ws.Range("E48").Value = '  -9.15%  '
$ws.Range("E49").Value = '  -6.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.77'
$ws.Range("E51").Value = '  -9.92%  '
